$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (row 7) describing a new logged test-mail entry, mirroring
# the structure of the existing rows (e.g. row 6).

$newAntwoord = "Beste klant,`n" +
    "Bedankt voor je bericht. Ik begrijp dat je retour nog niet is verwerkt en ik wil je graag helpen om dit op te lossen.`n" +
    "Om dit verder te onderzoeken, heb ik wat meer informatie nodig. Zou je zo vriendelijk willen zijn om je ordernummer en/of trackingnummer met me te delen? Hiermee kan ik de status van je retour nakijken en je zo goed mogelijk van dienst zijn.`n" +
    "Ik kijk uit naar je reactie.`n" +
    "Met vriendelijke groet,`n" +
    "[Naam van de e-mailassistent]  `n" +
    "Jamie  `n" +
    "Nederlandse e-mailassistent  `n" +
    "[Bedrijfsnaam]"

$ws.Range("A7").Value = "Testmail #11: Mijn retour is nog steeds niet verwerkt."
$ws.Range("B7").Value = $newAntwoord
$ws.Range("C7").Value = "Mijn retour is nog steeds niet verwerkt."
$ws.Range("D7").Value = "mailmind.test@zohomail.eu"
$ws.Range("E7").Value = "Retour / Terugbetaling"
$ws.Range("F7").Value = "2025-08-05 18:28:24"
$ws.Range("G7").Value = "Ja"
$ws.Range("H7").Value = "Nee"
$ws.Range("I7").Value = "Ja"
$ws.Range("J7").Value = "Nee"

# The multi-line "Antwoord" text otherwise leaves the row with an
# auto-expanded custom height; re-fit it back to the sheet default so the
# row matches the compact layout used by the other rows.
$ws.Rows.Item(7).AutoFit()
